# Update "想去人数" (interest count) values in the F column across the
# three sheets that list events: 展览 (sheet1), 本地生活 (sheet3) and the
# combined 全部类型 (sheet4). 演出 (sheet2) is untouched.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws3 = $wb.Worksheets.Item("本地生活")
$ws4 = $wb.Worksheets.Item("全部类型")

# --- 展览 sheet ---
$ws1.Range("F2").Value = 257
$ws1.Range("F3").Value = 2713
$ws1.Range("F7").Value = 2292
$ws1.Range("F8").Value = 1837
$ws1.Range("F9").Value = 219
$ws1.Range("F11").Value = 2483
$ws1.Range("F12").Value = 553
$ws1.Range("F15").Value = 31
$ws1.Range("F17").Value = 116
$ws1.Range("F18").Value = 9234
$ws1.Range("F20").Value = 7173
$ws1.Range("F21").Value = 11727
$ws1.Range("F24").Value = 233
$ws1.Range("F25").Value = 356
$ws1.Range("F26").Value = 561
$ws1.Range("F27").Value = 2595
$ws1.Range("F28").Value = 234
$ws1.Range("F30").Value = 2547
$ws1.Range("F31").Value = 713
$ws1.Range("F33").Value = 4513
$ws1.Range("F34").Value = 914
$ws1.Range("F35").Value = 355
$ws1.Range("F37").Value = 529

# --- 本地生活 sheet ---
$ws3.Range("F2").Value = 629
$ws3.Range("F4").Value = 156

# --- 全部类型 sheet ---
$ws4.Range("F2").Value = 629
$ws4.Range("F3").Value = 257
$ws4.Range("F5").Value = 2713
$ws4.Range("F9").Value = 2292
$ws4.Range("F11").Value = 1837
$ws4.Range("F13").Value = 219
$ws4.Range("F14").Value = 2483
$ws4.Range("F16").Value = 553
$ws4.Range("F21").Value = 116
$ws4.Range("F22").Value = 9234
$ws4.Range("F24").Value = 7174
$ws4.Range("F25").Value = 11727
$ws4.Range("F28").Value = 233
$ws4.Range("F29").Value = 356
$ws4.Range("F31").Value = 561
$ws4.Range("F33").Value = 2595
$ws4.Range("F36").Value = 234
$ws4.Range("F39").Value = 4513
$ws4.Range("F46").Value = 529
